# Edit script: add new Japanese names to row 13, add new Hebrew row 17,
# add conditional formatting (duplicate-values highlighting) for the new rows,
# widen the default column width, and update the view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Append 72 new Japanese names to row 13 (columns AK:DD) ---
$japaneseNew = @("Katsushika","Ken","Nijiro","Yoshiki","Hayao","Shigeru","Yoko","Akihito","Hirohito","Akira","Yoshi","Shinz","Kae","Kento","Kendo","Hiroshi","Korin","Takashi","Kakashi","Yoshimoto","Kazuo","Shinro","Hiroaki","Toshio","Monzo","Masafumi","Kosaku","Boshiro","Enomoto","Inoue","Hikojiro","Takatsugu","Kabayama","Kakuji","Kato","Keiji","Mineichi","Takeo","Ryu","Ryunosuke","Gunichi","Chuichi","Shoji","Kondo","Jisaku","Keisuke","Takijiro","Jisaburo","Sanemi","Inosuke","Giyu","Zenitsu","Tokito","Akaza","Genya","Muzan","Kyojuro","Obanai","Sakonji","Neji","Shikamaru","Nagato","Obito","Dan","Kiba","Miyamoto","Honda","Shinmen","Musashi","Shiro","Taira","Utagawa")
$col = 37  # column AK
foreach ($v in $japaneseNew) {
    $ws.Cells.Item(13, $col).Value = $v
    $col = $col + 1
}

# --- 2) Add new row 17 (hebrew category, columns A:DB) ---
$hebrewRow = @("hebrew","Adan","Amos","Benjamin","Esdras","Jacob","Levi","Moshe","Natan","Barack","Biel","Cain","Caleb","Carmelo","Zev","Ariel","Baruc","Abraham","Isaac","Jose","David","Jesus","Juan","Mateo","Jeremias","Samuel","Ezequiel","Pablo","Bernabe","Timoteo","Lucas","Abba","Adam","Aharon","Akiva","Alexander","Pedro","Alon","Amram","Arie","Avi","Avigdor","Avner","Azriel","Baruj","Ben Tzion","Berel","Betzalel","Boaz","Calev","Carmi","Dan","Doron","Daniel","Ehud","Dov","Eitan","Elazar","Eljanan","Eldad","Elimelej","Elisha","Eliakim","Emanuel","Ezra","Faivel","Fishel","Gad","Gamliel","Gabriel","Gershon","Gidon","Hillel","Hirch","Lerajmiel","Ilan","Josef","Itsjak","Israel","Issur","Itamar","Janan","Jonathan","Janoj","Kalman","Kalonimos","Leib","Lior","Mendel","Menashe","Mijael","Mordejai","Natanel","Nejemia","Rafael","Seth","Shabtai","Shaul","Shalom","Shimon","Tamir","Uriel","Uziel","Yair","Yoel","Zalman")
$col = 1
foreach ($v in $hebrewRow) {
    $ws.Cells.Item(17, $col).Value = $v
    $col = $col + 1
}

# --- 3) Conditional formatting: duplicate-value highlighting for new rows ---
$fc13 = $ws.Range("A13:XFD13").FormatConditions.AddUniqueValues()
$fc13.DupeUnique = 1
$fc13.Interior.Color = 39423
$fc13.Priority = 5

$fc17 = $ws.Range("A17:XFD17").FormatConditions.AddUniqueValues()
$fc17.DupeUnique = 1
$fc17.Interior.Color = 39423
$fc17.Priority = 4

# --- 4) Widen default column width across the sheet ---
$ws.Columns.ColumnWidth = 15.8888888888889
$ws.StandardWidth = 15.8888888888889

# --- 5) Update view/selection state (scrolled + selected cell) ---
$null = $ws.Range("DD13").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 101
$win.ScrollRow = 4

Write-Host "edit complete"
